# Update "想去人数" (number of people interested) figures for several
# conventions/events, in both the "展览" sheet and the aggregated
# "全部类型" sheet, reflecting the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) -------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F7").Value  = 4581
$ws1.Range("F10").Value = 103
$ws1.Range("F12").Value = 87
$ws1.Range("F13").Value = 699
$ws1.Range("F14").Value = 186
$ws1.Range("F15").Value = 992
$ws1.Range("F20").Value = 118
$ws1.Range("F21").Value = 94
$ws1.Range("F22").Value = 3525
$ws1.Range("F23").Value = 5884
$ws1.Range("F29").Value = 3362
$ws1.Range("F30").Value = 362
$ws1.Range("F36").Value = 215
$ws1.Range("F37").Value = 263
$ws1.Range("F39").Value = 127
$ws1.Range("F40").Value = 1015
$ws1.Range("F41").Value = 912
$ws1.Range("F42").Value = 20
$ws1.Range("F44").Value = 48
$ws1.Range("F45").Value = 51
$ws1.Range("F46").Value = 470

# --- Sheet "全部类型" (All Types) -----------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F7").Value  = 4581
$ws4.Range("F10").Value = 103
$ws4.Range("F13").Value = 87
$ws4.Range("F14").Value = 699
$ws4.Range("F15").Value = 186
$ws4.Range("F16").Value = 992
$ws4.Range("F21").Value = 118
$ws4.Range("F22").Value = 94
$ws4.Range("F23").Value = 3525
$ws4.Range("F24").Value = 5884
$ws4.Range("F30").Value = 3362
$ws4.Range("F31").Value = 362
$ws4.Range("F37").Value = 215
$ws4.Range("F38").Value = 263
$ws4.Range("F40").Value = 127
$ws4.Range("F41").Value = 1015
$ws4.Range("F42").Value = 912
$ws4.Range("F43").Value = 20
$ws4.Range("F45").Value = 48
$ws4.Range("F46").Value = 51
$ws4.Range("F47").Value = 470

$wb.Save()
